# Automatically re-generate list and index
# Shift every "Review date" value that falls on the 27th of its month
# forward by one day (the 27th -> the 28th), leaving every other
# review date untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $current = [string]$cell.Value2

    if ($current -match '^(\d{4}-\d{2})-27$') {
        $newValue = $Matches[1] + "-28"

        # Force the cell to keep storing a plain text value (instead of
        # Excel auto-converting the "yyyy-mm-dd" string into a real date
        # serial number), then restore the cell's original style so no
        # visible formatting changes are introduced.
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.Style = $origStyle
    }
}
